# added 4wk low sales check
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row -> [D (MyForecast), H (Inventory Coverage), L (Seasonality Index)]
$forecastRows = @{
    2  = @(362, 14.52, 1.2)
    3  = @(299, 16.37, 0.98)
    4  = @(274, 16.78, 0.9)
    5  = @(292, 14.8, 1.04)
    6  = @(291, 13.85, 0.99)
    7  = @(263, 14.22, 1.03)
    8  = @(225, 15.45, 1.13)
    9  = @(212, 15.34, 0.83)
    10 = @(245, 12.41, 1.15)
    11 = @(289, 9.67, 0.82)
    12 = @(300, 8.35, 1.05)
    13 = @(295, 7.48, 0.91)
    14 = @(273, 7, 1.16)
    15 = @(270, 6.07, 1.09)
    16 = @(286, 4.78, 1.2)
    17 = @(320, 3.38, 0.87)
}

foreach ($row in $forecastRows.Keys) {
    $vals = $forecastRows[$row]
    $ws1.Cells.Item($row, 4).Value = $vals[0]   # D = MyForecast
    $ws1.Cells.Item($row, 8).Value = $vals[1]   # H = Inventory Coverage
    $ws1.Cells.Item($row, 12).Value = $vals[2]  # L = Seasonality Index
}

# --- Sheet 2: "Summary" ---
# Column B on this sheet stores metrics as TEXT (t="inlineStr"/shared string),
# even though the text itself looks numeric. Force text entry (like typing
# into a cell pre-formatted as Text) so the written value keeps its string
# type instead of being auto-coerced to a number, then restore the default
# "Normal" style so no stray number-format residue is left on the cell.
$ws2 = $wb.Worksheets.Item("Summary")

$summaryRows = @{
    9  = "4496"  # Total Forecast (16 Weeks)
    10 = "2218"  # Total Forecast (8 Weeks)
    11 = "1227"  # Total Forecast (4 Weeks)
    12 = "362"   # Max Forecast
    14 = "212"   # Min Forecast
}

foreach ($row in $summaryRows.Keys) {
    $cell = $ws2.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryRows[$row]
    $cell.Style = "Normal"
}
